$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final player/position/team table (rows 2-19), reflecting the roster
# update: two new players added (Spencer Dinwiddie, Davion Mitchell) right
# after Shai Gilgeous-Alexander, and Bam Adebayo / Kyrie Irving / Jordan
# Clarkson moved down to the bottom of the list.
$data = @(
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Shai Gilgeous-Alexander", "PG,SG", "Oklahoma City Thunder"),
    @("Spencer Dinwiddie", "PG,SG", "Dallas Mavericks"),
    @("Davion Mitchell", "PG,SG", "Miami Heat"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Zach LaVine", "SG,SF", "Sacramento Kings"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Rui Hachimura", "SF,PF", "Los Angeles Lakers"),
    @("Kyle Kuzma", "PF", "Milwaukee Bucks"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Gradey Dick", "SG,SF", "Toronto Raptors"),
    @("Matas Buzelis", "SF,PF", "Chicago Bulls"),
    @("Joel Embiid", "C", "Philadelphia 76ers"),
    @("Bam Adebayo", "PF,C", "Miami Heat"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Jordan Clarkson", "SG,SF", "Utah Jazz")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
